$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '91.777.52'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  +0.91%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '3.104.67'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('E4').Value = '  -0.02%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '235.70'
$cell.Style = 'Normal'
$ws.Range('E5').Value = '  -1.67%  '
$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '614.12'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.74%  '
$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '1.10'
$cell.Style = 'Normal'
$ws.Range('E7').Value = '  -2.21%  '
$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.387'
$cell.Style = 'Normal'
$ws.Range('E8').Value = '  +3.43%  '
$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.999'
$cell.Style = 'Normal'
$ws.Range('E9').Value = '  -0.08%  '
$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '3.097.32'
$cell.Style = 'Normal'
$ws.Range('E10').Value = '  -1.64%  '
$cell = $ws.Range('D11')
$cell.NumberFormat = '@'
$cell.Value = '0.733'
$cell.Style = 'Normal'
$ws.Range('E11').Value = '  -1.64%  '
$ws.Range('E12').Value = '  -1.41%  '
$ws.Range('E13').Value = '  -0.42%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '91.977.22'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  +1.16%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '33.91'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -3.50%  '
$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '5.42'
$cell.Style = 'Normal'
$ws.Range('E16').Value = '  -2.61%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '3.690.67'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -1.48%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '3.086.42'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  -2.10%  '
$cell = $ws.Range('D19')
$cell.NumberFormat = '@'
$cell.Value = '3.73'
$cell.Style = 'Normal'
$ws.Range('E19').Value = '  -0.27%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '14.56'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -4.02%  '
$ws.Range('E21').Value = '  -4.82%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '9.32'
$cell.Style = 'Normal'
$ws.Range('E22').Value = '  +1.87%  '
$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$cell = $ws.Range('D23')
$cell.NumberFormat = '@'
$cell.Value = '444.09'
$cell.Style = 'Normal'
$ws.Range('E23').Value = '  -1.39%  '
$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '0.0000194'
$cell.Style = 'Normal'
$ws.Range('E24').Value = '  -5.05%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '5.72'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  -0.45%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '86.01'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  -3.12%  '
$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '11.63'
$cell.Style = 'Normal'
$ws.Range('E27').Value = '  -3.14%  '
$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '3.261.30'
$cell.Style = 'Normal'
$ws.Range('E28').Value = '  -1.91%  '
$ws.Range('E29').Value = '  -0.01%  '
$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '0.133'
$cell.Style = 'Normal'
$ws.Range('E30').Value = '  -6.07%  '
$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '0.232'
$cell.Style = 'Normal'
$ws.Range('E31').Value = '  -1.63%  '
$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '0.168'
$cell.Style = 'Normal'
$ws.Range('E32').Value = '  -1.86%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '9.08'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -3.11%  '
$ws.Range('E34').Value = '  +2.21%  '
$ws.Range('E35').Value = '  -8.65%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '26.02'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -1.40%  '
$ws.Range('E37').Value = '  +1.32%  '
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '1.89'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -4.31%  '
$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '480.31'
$cell.Style = 'Normal'
$ws.Range('E39').Value = '  -5.53%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '1.29'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -3.73%  '
$ws.Range('B41').Value = 'WhiteBITCoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '23.85'
$cell.Style = 'Normal'
$ws.Range('E41').Value = '  +7.95%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.430'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -3.91%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '3.26'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -5.49%  '
$ws.Range('B44').Value = 'Binance-PegBSC-USD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '0.758'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -24.29%  '
$ws.Range('E45').Value = '  +0.04%  '
$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '164.04'
$cell.Style = 'Normal'
$ws.Range('E46').Value = '  +5.36%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.687'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  -4.20%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '1.88'
$cell.Style = 'Normal'
$ws.Range('E48').Value = '  -3.15%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '1.38'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  +0.47%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '0.0334'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  +2.93%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '4.39'
$cell.Style = 'Normal'
$ws.Range('E51').Value = '  -1.75%  '
